# Skill system / PlayerCardPage
# Adds a new "系列最大點數" (max series points) column (L) to the Skill sheet,
# populates it for every skill row, and sets the sheet page setup to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell L1, styled like the other header cells (K1 etc. use the
# 游ゴシック header font).
$ws.Range("L1").Value = "系列最大點數"
$ws.Range("L1").Font.Name = "游ゴシック"
$ws.Range("L1").Font.Charset = 134

# Per-row values for the new column (rows 2-31 correspond to the skill rows).
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 3
    8  = 2
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 3
    20 = 3
    21 = 2
    22 = 1
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 12).Value = $values[$row]
}

# Select L19 like the saved workbook did.
$ws.Range("L19").Select()

# Page setup: portrait orientation.
$ws.PageSetup.Orientation = 1
